$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
$ws.Cells.Item(1,2).Value = "Simplified Text"
$ws.Cells.Item(1,4).Value = "Readability Simplified"
$ws.Cells.Item(1,5).Value = "Words Original"
$ws.Cells.Item(1,6).Value = "Words Simplified"
$ws.Cells.Item(1,7).Value = "Sentences Original"
$ws.Cells.Item(1,8).Value = "Sentences Simplified"

# Copy the header style (bold, bordered, centered) used by A1:D1 onto the new E1:H1 headers
$ws.Range("A1").Copy()
$ws.Range("E1:H1").PasteSpecial(-4122)

# --- New "Words Original", "Words Simplified", "Sentences Original", "Sentences Simplified" data columns (E:H), rows 2-100 ---
$data = @(
    @(2,649,281,43,31),
    @(3,827,165,41,16),
    @(4,835,182,38,14),
    @(5,791,290,42,20),
    @(6,554,176,30,11),
    @(7,256,71,14,6),
    @(8,475,224,26,14),
    @(9,191,115,10,7),
    @(10,628,191,32,13),
    @(11,151,109,9,8),
    @(12,780,318,40,16),
    @(13,552,268,29,15),
    @(14,514,286,28,18),
    @(15,254,207,16,14),
    @(16,87,58,6,5),
    @(17,626,169,42,13),
    @(18,849,205,41,17),
    @(19,503,178,30,14),
    @(20,292,137,13,12),
    @(21,297,212,14,11),
    @(22,507,173,22,15),
    @(23,549,158,30,12),
    @(24,580,267,23,15),
    @(25,507,214,22,13),
    @(26,828,132,49,11),
    @(27,791,316,39,21),
    @(28,265,69,17,5),
    @(29,455,273,23,18),
    @(30,402,288,15,15),
    @(31,95,73,8,8),
    @(32,854,240,41,19),
    @(33,283,187,18,13),
    @(34,467,173,25,11),
    @(35,344,114,15,10),
    @(36,295,141,15,10),
    @(37,742,173,37,13),
    @(38,569,320,27,18),
    @(39,479,107,25,10),
    @(40,515,279,24,15),
    @(41,736,119,42,9),
    @(42,539,213,24,14),
    @(43,824,296,45,15),
    @(44,801,134,36,10),
    @(45,282,107,17,13),
    @(46,539,384,29,23),
    @(47,804,255,44,17),
    @(48,530,277,26,18),
    @(49,461,233,17,13),
    @(50,572,434,24,23),
    @(51,742,289,33,18),
    @(52,743,286,39,17),
    @(53,517,274,22,19),
    @(54,604,267,24,15),
    @(55,840,257,34,14),
    @(56,575,304,31,21),
    @(57,817,172,35,23),
    @(58,719,273,39,21),
    @(59,133,142,9,9),
    @(60,519,214,29,13),
    @(61,768,186,101,16),
    @(62,47,30,3,3),
    @(63,769,164,52,13),
    @(64,599,206,26,14),
    @(65,462,247,20,14),
    @(66,807,206,56,13),
    @(67,838,263,44,15),
    @(68,131,106,8,8),
    @(69,796,243,47,17),
    @(70,808,251,39,15),
    @(71,591,323,30,21),
    @(72,587,326,39,25),
    @(73,730,235,41,14),
    @(74,330,156,15,10),
    @(75,427,195,24,14),
    @(76,700,246,32,15),
    @(77,614,266,25,16),
    @(78,826,359,45,20),
    @(79,885,200,38,11),
    @(80,373,225,25,16),
    @(81,642,207,23,14),
    @(82,517,140,23,11),
    @(83,244,182,15,13),
    @(84,800,127,58,10),
    @(85,753,114,29,9),
    @(86,620,190,38,14),
    @(87,786,152,36,10),
    @(88,818,150,39,12),
    @(89,640,164,29,12),
    @(90,550,243,29,16),
    @(91,641,235,25,22),
    @(92,550,259,26,15),
    @(93,672,145,29,8),
    @(94,819,212,41,13),
    @(95,631,200,31,18),
    @(96,555,100,25,8),
    @(97,526,96,31,6),
    @(98,768,219,39,12),
    @(99,673,288,29,18),
    @(100,821,307,53,17)
)
foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 5).Value = $row[1]
    $ws.Cells.Item($r, 6).Value = $row[2]
    $ws.Cells.Item($r, 7).Value = $row[3]
    $ws.Cells.Item($r, 8).Value = $row[4]
}

# --- Strip leading "Corrected text:" / "Simplified text:" / "Subject: " labels that prefixed some entries ---
$a35 = $ws.Cells.Item(35,1).Text
$prefixA35 = "Corrected text:`r`n`r`n"
if ($a35.StartsWith($prefixA35)) {
    $ws.Cells.Item(35,1).Value = $a35.Substring($prefixA35.Length)
}

$b35 = $ws.Cells.Item(35,2).Text
$prefixB35 = "Simplified text:`r`n`r`n"
if ($b35.StartsWith($prefixB35)) {
    $ws.Cells.Item(35,2).Value = $b35.Substring($prefixB35.Length)
}

$b42 = $ws.Cells.Item(42,2).Text
$prefixB42 = "Subject: "
if ($b42.StartsWith($prefixB42)) {
    $ws.Cells.Item(42,2).Value = $b42.Substring($prefixB42.Length)
}

$a45 = $ws.Cells.Item(45,1).Text
$prefixA45 = "Corrected text:`r`n`r`n"
if ($a45.StartsWith($prefixA45)) {
    $ws.Cells.Item(45,1).Value = $a45.Substring($prefixA45.Length)
}

# --- Recalculated readability scores for the rows whose text changed ---
$ws.Cells.Item(35,3).Value = 16.54237209302326
$ws.Cells.Item(35,4).Value = 11.42091228070176
$ws.Cells.Item(42,4).Value = 15.55014419852448
$ws.Cells.Item(45,3).Value = 16.11132665832291

# --- Selection / view bookkeeping to match the saved state ---
$ws.Range("A1").Select()
